$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F9").Value = 101787
$ws.Range("F10").Value = 90267
